$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates from the latest cryptos data refresh.
# Cells whose new text looks like a plain number (e.g. "1.00") are forced
# to Text format first so Excel keeps them as strings (matching the source
# data, which stores prices as text) instead of silently converting them to
# numeric values; the temporary format is cleared again right after so no
# extra style is left behind on the cell.

$ws.Range('D2').Value = '56.422.78'
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').Value = '2.980.81'
$ws.Range('E3').Value = '  -2.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '502.47'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.65'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.86%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('E9').Value = '  +3.38%  '
$ws.Range('E10').Value = '  +2.13%  '
$ws.Range('E11').Value = '  -2.36%  '
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').Value = '3.488.15'
$ws.Range('E13').Value = '  -3.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.07'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.71%  '
$ws.Range('D15').Value = '56.411.90'
$ws.Range('E15').Value = '  +3.18%  '
$ws.Range('E16').Value = '  +3.82%  '
$ws.Range('D17').Value = '2.977.52'
$ws.Range('E17').Value = '  -3.63%  '
$ws.Range('E18').Value = '  +3.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.35'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('E20').Value = '  +2.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '325.40'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -3.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.11'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -5.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.163'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.39%  '
$ws.Range('D27').Value = '0.0₃0889'
$ws.Range('E27').Value = '  +1.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.46'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.78'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.08%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.17'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.55%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.74'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.24'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '156.06'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  -2.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.27'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.52'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0671'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.10'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('D40').Value = '3.014.44'
$ws.Range('E40').Value = '  -3.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '35.99'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('E43').Value = '  -3.32%  '
$ws.Range('D44').Value = '2.236.40'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.986'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.70%  '
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('E47').Value = '  -4.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.93'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +12.85%  '
$ws.Range('E49').Value = '  +4.09%  '
$ws.Range('E50').Value = '  -2.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.89'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.12%  '
